# Update feature_codebook sheet to reflect the "drabo simulation" codebook.
# The table is rebuilt row-by-row (row 2 stays blank, exactly like before).
# Column A/B only carry a value on the first row of each variable's block;
# C/D repeat for every category row of a categorical/boolean variable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @(A, B, C, D)  ($null means "leave/clear that cell")
$rows = @{
    1  = @("variable_name", "type", "values", "description")
    3  = @("age_min", "integer", "0-100", "minimum age of demographic group")
    4  = @("age_max", "integer", "0-100", "maximum age of demographic group")
    5  = @("race", "categorical", 0, "white")
    6  = @($null, $null, 1, "black")
    7  = @($null, $null, 2, "latinX")
    8  = @("msm", "boolean", "0/1", "men who have sex with men")
    9  = @("alive", "boolean", "0/1", "alive in the current time period")
    10 = @("hiv", "categorical", 0, "uninfected")
    11 = @($null, $null, 1, "primary")
    12 = @($null, $null, 2, "asymptomatic")
    13 = @($null, $null, 3, "symptomatic")
    14 = @($null, $null, 4, "AIDS")
    15 = @("aware", "boolean", "0/1", "aware of serostatus")
    16 = @("prep", "boolean", 0, "not on PrEP (pre-exposure prophylaxis)")
    17 = @($null, $null, 1, "on PrEP")
    18 = @("art", "categorical", 0, "not on ART (anti-retroviral therapy)")
    19 = @($null, $null, 1, "on ART")
}

# The sheet used to extend to row 20; clear everything below the new last
# row (19) first so no stale values/old dimension linger.
$ws.Range("A2:D20").ClearContents()

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($c = 0; $c -lt 4; $c++) {
        $v = $vals[$c]
        $cell = $ws.Cells.Item($r, $c + 1)
        if ($null -eq $v) {
            $cell.ClearContents()
        } else {
            $cell.Value = $v
        }
    }
}

# View/selection cosmetics captured in the diff.
$ws.Cells.Item(1, 1).Select() | Out-Null
$ws.Range("B1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 200
$ws.Range("D20").Select() | Out-Null
